# Update workbook/worksheet for data through 2022-03-06 (adding 2022-03-14 carjacking data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2022-03-06"

# Update the header label in B1 to match the new "through" date
$ws.Range("B1").Value = "March 2022 (through March 06)"

# Apply the numeric cell updates for the neighborhoods affected by the new day's data
$ws.Range("E4").Value  = 2   # North Lawndale

$ws.Range("B5").Value  = 2   # Garfield Park
$ws.Range("H5").Value  = 1
$ws.Range("N5").Value  = 3
$ws.Range("Q5").Value  = 3

$ws.Range("B6").Value  = 1   # Rogers Park

$ws.Range("B7").Value  = 1   # South Shore

$ws.Range("H10").Value = 1   # Chicago Lawn

$ws.Range("B12").Value = 3   # Englewood
$ws.Range("T12").Value = 2

$ws.Range("B13").Value = 1   # Woodlawn

$ws.Range("Q15").Value = 1   # Humboldt Park

$ws.Range("W16").Value = 1   # Little Italy, UIC

$ws.Range("H18").Value = 2   # Washington Heights

$ws.Range("H24").Value = 1   # Wicker Park

$ws.Range("Q33").Value = 1   # Belmont Cragin

$ws.Range("H44").Value = 2   # Grand Boulevard

$ws.Range("E55").Value = 1   # Clearing

$ws.Range("Q82").Value = 1   # South Chicago
